# "add M45 and M102" -- insert two missing Messier catalog rows into the
# Messier sheet, shifting subsequent rows down, and update the saved
# view state (active sheet / selection) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Messier")

# --- Insert the M45 row at sheet row 44 (pushes the old row 44.. data down by one) ---
$ws.Rows.Item(44).Insert()

# --- Insert the M102 row at sheet row 101 (post-shift position) ---
$ws.Rows.Item(101).Insert()

# Fill in new row 44 = M45 entry (column A intentionally left blank, as in the
# surrounding data the new row has no NGC cross-reference yet)
$ws.Cells.Item(44, 2).Value = 1
$ws.Cells.Item(44, 3).Value = 45
$ws.Cells.Item(44, 4).Value = 160
$ws.Cells.Item(44, 5).Value = 3.783083
$ws.Cells.Item(44, 6).Value = 24.1144
$ws.Cells.Item(44, 7).Formula = '=RIGHT("              "&TRUNC((D44+200)/10),4)'
$ws.Cells.Item(44, 8).Formula = '=RIGHT("              " &TRUNC(E44*65536 / 24),6)'
$ws.Cells.Item(44, 9).Formula = '=RIGHT("              " & TRUNC(F44*32767 / 90),6)'

# Fill in new row 101 = M102 entry
$ws.Cells.Item(101, 2).Value = 0
$ws.Cells.Item(101, 3).Value = 102
$ws.Cells.Item(101, 4).Value = 989
$ws.Cells.Item(101, 5).Value = 15.1079583
$ws.Cells.Item(101, 6).Value = 55.76555555
$ws.Cells.Item(101, 7).Formula = '=RIGHT("              "&TRUNC((D101+200)/10),4)'
$ws.Cells.Item(101, 8).Formula = '=RIGHT("              " &TRUNC(E101*65536 / 24),6)'
$ws.Cells.Item(101, 9).Formula = '=RIGHT("              " & TRUNC(F101*32767 / 90),6)'

# --- Update the saved view/selection state ---
$ws.Activate()
$ws.Range("G103").Select() | Out-Null
